# feat: integração e verificador de cartões
#
# 1) "Cartão verifier": Guilherme Alves' CPF (row 8) was a bogus/duplicate
#    value ("369.852.147-00"); fix it to the correct CPF shared with
#    Felipe Mendes ("852.963.741-28"). That retires the old unique string
#    from the shared-strings table entirely.
# 2) "Integração": append a new record (row 12) for Juliana Costa at the
#    BTG brokerage, re-using her existing CPF/status/profile/planner but
#    with a freshly generated card/CPF number for the new integration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the CPF on row 8 (Guilherme Alves) ---
$ws.Range("C8").Value = "852.963.741-28"

# --- 2. Add a new row 12, cloning row 11's formatting then editing cells ---
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A12").Value = "Juliana Costa"
$ws.Range("B12").Value = "BTG"
$ws.Range("C12").Value = "987.654.321-01"
$ws.Range("D12").Value = "Portabilidade"
$ws.Range("E12").Value = "Moderado"
$ws.Range("F12").Value = "Gabrihel Bieguelman"

# --- 3. Match the last-saved selection recorded in the sheet view ---
$ws.Range("C8").Select() | Out-Null

Write-Output "applied integration + card verification edits"
